$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thomas")

# New "Gewerkt" label next to the Study-row total formula
$ws.Range("F3").Value = "Gewerkt"

# Expand the running-total formula in F18 to cover the newly added rows,
# and add the "Gedeclareerd" label next to it
$ws.Range("F18").Formula = "=SUM(C2:C21)"
$ws.Range("G18").Value = "Gedeclareerd"

# New row 26
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("A26").Value = "Meeting data formatting MSC Maersk Routescanner"
$ws.Range("B26").Value = (Get-Date -Year 2023 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C26").Value = 3

# New row 27
$ws.Range("B25").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Range("A27").Value = "Combining dataframes in python"
$ws.Range("B27").Value = (Get-Date -Year 2023 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C27").Value = 2.5

# Update the saved selection state to match the authored session
$ws.Range("E4").Select() | Out-Null
